$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.132.49'
$ws.Range('E2').Value = '  +3.27%  '
$ws.Range('D3').Value = '2.302.48'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.86%  '
$ws.Range('E7').Value = '  +2.73%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.522'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.65%  '
$ws.Range('E10').Value = '  +5.14%  '
$ws.Range('E11').Value = '  +5.41%  '
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.12'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.31%  '
$ws.Range('D14').Value = '2.657.84'
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.72%  '
$ws.Range('D16').Value = '2.300.06'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('E17').Value = '  +3.13%  '
$ws.Range('D18').Value = '43.044.30'
$ws.Range('E18').Value = '  +3.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.52'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.41%  '
$ws.Range('E20').Value = '  +3.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('E24').Value = '  +5.72%  '
$ws.Range('E25').Value = '  +2.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.67%  '
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.35%  '
$ws.Range('E32').Value = '  +3.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.80'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0740'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.107'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.80%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.39'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.87%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.116'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.29'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.22%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.34'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.54%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.971.81'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0289'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('E46').Value = '  +4.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.76'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.19%  '
$ws.Range('E49').Value = '  +16.15%  '
$ws.Range('D50').Value = '2.526.64'
$ws.Range('E50').Value = '  +1.65%  '
$ws.Range('E51').Value = '  +3.24%  '
